# Update the BLS Data Series sheet: replace the "Series ID" column header
# and the raw BLS series-ID row labels with human readable demographic
# category labels, for the Project 2 proportions/prediction plot update.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BLS Data Series")

# Header
$ws.Range("A1").Value = "Category"

# Row labels (BLS series IDs -> demographic categories)
$ws.Range("A2").Value = "White Men"
$ws.Range("A3").Value = "White Women"
$ws.Range("A5").Value = "Asian Women"
$ws.Range("A6").Value = "Hispanic Women"

# "Black Women" was retyped without carrying over the shaded/Arial row-label
# style, so it now renders in the workbook's default bold font with no fill.
$ws.Range("A4").Value = "Black Women"
$ws.Range("A4").ClearFormats()
$ws.Range("A4").Font.Bold = $true

# Restore the active selection to A2 (as captured in the saved workbook view)
$ws.Range("A2").Select() | Out-Null
